$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date on Overview sheet
$wsOverview.Range("G2").Value = "2016-10-27 09:56:27"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-10-27 09:56:14"
$wsZhCn.Range("K2").Value = "2016-10-27 09:56:53"

# de-de sheet: Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-10-27 09:57:10"
